$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 5 (old rows 5-17 shift down to 7-19)
$ws.Rows("5:6").Insert()

# Fill in the new "comision 1" / "comision 2" rows
$ws.Range("F5").Value = "comision 1"
$ws.Range("G5").Value = 157.5
$ws.Range("H5").Value = 1
$ws.Range("I5").Formula = "=G5*H5"

$ws.Range("F6").Value = "comision 2"
$ws.Range("G6").Value = 337.5
$ws.Range("H6").Value = 1
$ws.Range("I6").Formula = "=G6*H6"

# Update the selection to match the target state
$ws.Range("G5:G6").Select()
